$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 6946
$ws.Range("J10").Value = 6946
$ws.Range("L10").Value = 6946
$ws.Range("N10").Value = -7532
$ws.Range("H19").Value = 1441.2142
$ws.Range("I19").Value = 1674.1818
$ws.Range("J19").Value = 587
$ws.Range("K19").Value = 1674.1818
$ws.Range("L19").Value = 587
$ws.Range("M19").Value = -1499.1818
$ws.Range("N19").Value = -937
$ws.Range("H31").Value = 276.8
$ws.Range("I31").Value = 276.8
$ws.Range("K31").Value = 830.4000000000001
$ws.Range("M31").Value = -600.4000000000001
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H137").Value = 2857.4285
$ws.Range("I137").Value = 2857.4285
$ws.Range("K137").Value = 8572.2855
$ws.Range("M137").Value = -6022.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 10000
$ws.Range("L56").Value = 10000
$ws.Range("N56").Value = -11484
$ws.Range("H68").Value = 29999
$ws.Range("J68").Value = 29999
$ws.Range("L68").Value = 29999
$ws.Range("N68").Value = -31621
$ws.Range("H71").Value = 29999
$ws.Range("J71").Value = 29999
$ws.Range("L71").Value = 89997
$ws.Range("N71").Value = -98109
$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -55059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 35000
$ws.Range("L21").Value = 35000
$ws.Range("N21").Value = -35472
$ws.Range("H22").Value = 399.5
$ws.Range("I22").Value = 297
$ws.Range("J22").Value = 502
$ws.Range("K22").Value = 297
$ws.Range("L22").Value = 502
$ws.Range("M22").Value = -124
$ws.Range("N22").Value = -848
$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -40990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 206.90909
$ws.Range("J22").Value = 214.33333
$ws.Range("L22").Value = 214.33333
$ws.Range("N22").Value = -914.3333299999999
$ws.Range("H25").Value = 15000
$ws.Range("I25").Value = 15000
$ws.Range("K25").Value = 15000
$ws.Range("M25").Value = -14826
$ws.Range("H97").Value = 35000
$ws.Range("J97").Value = 35000
$ws.Range("L97").Value = 35000
$ws.Range("N97").Value = -36982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 830.1579
$ws.Range("I2").Value = 408.33334
$ws.Range("J2").Value = 1553.2858
$ws.Range("K2").Value = 2450.00004
$ws.Range("L2").Value = 9319.714800000002
$ws.Range("M2").Value = -2337.00004
$ws.Range("N2").Value = -9545.714800000002
$ws.Range("H26").Value = 265.1111
$ws.Range("I26").Value = 160.5
$ws.Range("J26").Value = 1102
$ws.Range("K26").Value = 481.5
$ws.Range("L26").Value = 3306
$ws.Range("M26").Value = -193.5
$ws.Range("N26").Value = -3882
$ws.Range("H34").Value = 835.4545000000001
$ws.Range("I34").Value = 520
$ws.Range("J34").Value = 905.55554
$ws.Range("K34").Value = 1560
$ws.Range("L34").Value = 2716.66662
$ws.Range("M34").Value = -1476
$ws.Range("N34").Value = -2884.66662
$ws.Range("H40").Value = 142.875
$ws.Range("I40").Value = 161.25
$ws.Range("J40").Value = 124.5
$ws.Range("K40").Value = 645
$ws.Range("L40").Value = 498
$ws.Range("M40").Value = -576
$ws.Range("N40").Value = -636
$ws.Range("H81").Value = 1500
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 4500
$ws.Range("N81").Value = -6746
$ws.Range("H84").Value = 1500
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 13500
$ws.Range("N84").Value = -24732
$ws.Range("H103").Value = 4949.75
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 4949.75
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 14849.25
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -16607.25
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 5000
$ws.Range("J21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("N21").Value = -5348
$ws.Range("H22").Value = 1999.5
$ws.Range("I22").Value = 1999.5
$ws.Range("K22").Value = 1999.5
$ws.Range("M22").Value = -1704.5
$ws.Range("H27").Value = 1999.5
$ws.Range("I27").Value = 1999.5
$ws.Range("K27").Value = 1999.5
$ws.Range("M27").Value = -1892.5
$ws.Range("H30").Value = 850
$ws.Range("I30").Value = 850
$ws.Range("K30").Value = 850
$ws.Range("M30").Value = -742
$ws.Range("H92").Value = 35000
$ws.Range("J92").Value = 35000
$ws.Range("L92").Value = 35000
$ws.Range("N92").Value = -39992

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 22668.5
$ws.Range("J45").Value = 10626
$ws.Range("L45").Value = 10626
$ws.Range("N45").Value = -11608
$ws.Range("H81").Value = 5000
$ws.Range("I81").Value = 5000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 10000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -8939
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 5000
$ws.Range("I84").Value = 5000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 50000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -44696
$ws.Range("N84").ClearContents()
$ws.Range("H93").Value = 32500
$ws.Range("J93").Value = 32500
$ws.Range("L93").Value = 32500
$ws.Range("N93").Value = -37492
$ws.Range("H100").Value = 1572.5714
$ws.Range("I100").Value = 974
$ws.Range("K100").Value = 1948
$ws.Range("M100").Value = -1407
